# Applies the Titan_Profits leve-profit recompute (per scheduled runner diff).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 4195.3335
$ws.Cells.Item(6, 9).Value = 2812.2
$ws.Cells.Item(6, 10).Value = 11111
$ws.Cells.Item(6, 11).Value = 8436.599999999999
$ws.Cells.Item(6, 12).Value = 33333
$ws.Cells.Item(6, 13).Value = -8324.599999999999
$ws.Cells.Item(6, 14).Value = -33557

$ws.Cells.Item(15, 8).Value = 152482.7
$ws.Cells.Item(15, 9).Value = 152482.7
$ws.Cells.Item(15, 11).Value = 457448.1
$ws.Cells.Item(15, 13).Value = -457279.1

$ws.Cells.Item(17, 8).Value = 574735.7
$ws.Cells.Item(17, 9).Value = 250
$ws.Cells.Item(17, 10).Value = 643331
$ws.Cells.Item(17, 11).Value = 750
$ws.Cells.Item(17, 12).Value = 1929993
$ws.Cells.Item(17, 13).Value = -582
$ws.Cells.Item(17, 14).Value = -1930329

$ws.Cells.Item(21, 8).Value = 58000
$ws.Cells.Item(21, 10).Value = 58000
$ws.Cells.Item(21, 12).Value = 58000
$ws.Cells.Item(21, 14).Value = -58936

$ws.Cells.Item(23, 8).Value = 58000
$ws.Cells.Item(23, 10).Value = 58000
$ws.Cells.Item(23, 12).Value = 58000
$ws.Cells.Item(23, 14).Value = -58468

$ws.Cells.Item(32, 8).Value = 2704.7693
$ws.Cells.Item(32, 10).Value = 2704.7693
$ws.Cells.Item(32, 12).Value = 2704.7693
$ws.Cells.Item(32, 14).Value = -3356.7693

$ws.Cells.Item(55, 8).Value = 144.41667
$ws.Cells.Item(55, 9).Value = 171.125
$ws.Cells.Item(55, 10).Value = 91
$ws.Cells.Item(55, 11).Value = 171.125
$ws.Cells.Item(55, 12).Value = 91
$ws.Cells.Item(55, 13).Value = 42.875
$ws.Cells.Item(55, 14).Value = -519

$ws.Cells.Item(129, 8).Value = 1115.05
$ws.Cells.Item(129, 10).Value = 1183.9166
$ws.Cells.Item(129, 12).Value = 3551.7498
$ws.Cells.Item(129, 14).Value = -13551.7498

$ws.Cells.Item(137, 8).Value = 18519444
$ws.Cells.Item(137, 9).Value = 22222894
$ws.Cells.Item(137, 10).Value = 2197.2222
$ws.Cells.Item(137, 11).Value = 66668682
$ws.Cells.Item(137, 12).Value = 6591.6666
$ws.Cells.Item(137, 13).Value = -66666132
$ws.Cells.Item(137, 14).Value = -11691.6666

$ws.Cells.Item(138, 8).Value = 4306654.5
$ws.Cells.Item(138, 9).Value = 1198322.9
$ws.Cells.Item(138, 10).Value = 6668986.5
$ws.Cells.Item(138, 11).Value = 3594968.7
$ws.Cells.Item(138, 12).Value = 20006959.5
$ws.Cells.Item(138, 13).Value = -3589828.7
$ws.Cells.Item(138, 14).Value = -20017239.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(133, 8).Value = 52149.75
$ws.Cells.Item(133, 10).Value = 52149.75
$ws.Cells.Item(133, 12).Value = 52149.75
$ws.Cells.Item(133, 14).Value = -57209.75

$ws.Cells.Item(139, 8).Value = 39409.715
$ws.Cells.Item(139, 10).Value = 40894.668
$ws.Cells.Item(139, 12).Value = 40894.668
$ws.Cells.Item(139, 14).Value = -51174.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(16, 8).Value = 7000
$ws.Cells.Item(16, 10).Value = 10000
$ws.Cells.Item(16, 12).Value = 10000
$ws.Cells.Item(16, 14).Value = -10340

$ws.Cells.Item(59, 8).Value = 49999.5
$ws.Cells.Item(59, 10).Value = 49999.5
$ws.Cells.Item(59, 12).Value = 49999.5
$ws.Cells.Item(59, 14).Value = -51693.5

$ws.Cells.Item(80, 8).Value = 675.6
$ws.Cells.Item(80, 9).Value = 257.6
$ws.Cells.Item(80, 10).Value = 780.1
$ws.Cells.Item(80, 11).Value = 257.6
$ws.Cells.Item(80, 12).Value = 780.1
$ws.Cells.Item(80, 13).Value = 740.4
$ws.Cells.Item(80, 14).Value = -2776.1

$ws.Cells.Item(83, 8).Value = 675.6
$ws.Cells.Item(83, 9).Value = 257.6
$ws.Cells.Item(83, 10).Value = 780.1
$ws.Cells.Item(83, 11).Value = 1288
$ws.Cells.Item(83, 12).Value = 3900.5
$ws.Cells.Item(83, 13).Value = 3704
$ws.Cells.Item(83, 14).Value = -13884.5

$ws.Cells.Item(134, 8).Value = 20410514
$ws.Cells.Item(134, 9).Value = 23811338
$ws.Cells.Item(134, 11).Value = 71434014
$ws.Cells.Item(134, 13).Value = -71431479

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1994.4524
$ws.Cells.Item(31, 9).Value = 1170.3226
$ws.Cells.Item(31, 10).Value = 4317
$ws.Cells.Item(31, 11).Value = 1170.3226
$ws.Cells.Item(31, 12).Value = 4317
$ws.Cells.Item(31, 13).Value = -875.3226
$ws.Cells.Item(31, 14).Value = -4907

$ws.Cells.Item(34, 8).Value = 1994.4524
$ws.Cells.Item(34, 9).Value = 1170.3226
$ws.Cells.Item(34, 10).Value = 4317
$ws.Cells.Item(34, 11).Value = 1170.3226
$ws.Cells.Item(34, 12).Value = 4317
$ws.Cells.Item(34, 13).Value = -968.3226
$ws.Cells.Item(34, 14).Value = -4721

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 986.2093
$ws.Cells.Item(5, 9).Value = 589.1613
$ws.Cells.Item(5, 10).Value = 2011.9166
$ws.Cells.Item(5, 11).Value = 1767.4839
$ws.Cells.Item(5, 12).Value = 6035.7498
$ws.Cells.Item(5, 13).Value = -1655.4839
$ws.Cells.Item(5, 14).Value = -6259.7498

$ws.Cells.Item(12, 8).Value = 51.482758
$ws.Cells.Item(12, 9).Value = 64.375
$ws.Cells.Item(12, 10).Value = 35.615383
$ws.Cells.Item(12, 11).Value = 193.125
$ws.Cells.Item(12, 12).Value = 106.846149
$ws.Cells.Item(12, 13).Value = -20.125
$ws.Cells.Item(12, 14).Value = -452.846149

$ws.Cells.Item(38, 8).Value = 145.59091
$ws.Cells.Item(38, 9).Value = 100
$ws.Cells.Item(38, 10).Value = 162.6875
$ws.Cells.Item(38, 11).Value = 300
$ws.Cells.Item(38, 12).Value = 488.0625
$ws.Cells.Item(38, 13).Value = 47
$ws.Cells.Item(38, 14).Value = -1182.0625

$ws.Cells.Item(48, 8).Value = 2000
$ws.Cells.Item(48, 10).Value = 2000
$ws.Cells.Item(48, 12).Value = 6000
$ws.Cells.Item(48, 14).Value = -6500

$ws.Cells.Item(58, 8).Value = 7718.1177
$ws.Cells.Item(58, 9).Value = 1070
$ws.Cells.Item(58, 10).Value = 9142.714
$ws.Cells.Item(58, 11).Value = 3210
$ws.Cells.Item(58, 12).Value = 27428.142
$ws.Cells.Item(58, 13).Value = -3082
$ws.Cells.Item(58, 14).Value = -27684.142

$ws.Cells.Item(86, 8).Value = 2758.5
$ws.Cells.Item(86, 9).Value = 2758.5
$ws.Cells.Item(86, 10).Value = 0
$ws.Cells.Item(86, 11).Value = 8275.5
$ws.Cells.Item(86, 12).Value = 0
$ws.Cells.Item(86, 13).Value = -7089.5
$ws.Cells.Item(86, 14).ClearContents()

$ws.Cells.Item(89, 8).Value = 2758.5
$ws.Cells.Item(89, 9).Value = 2758.5
$ws.Cells.Item(89, 10).Value = 0
$ws.Cells.Item(89, 11).Value = 24826.5
$ws.Cells.Item(89, 12).Value = 0
$ws.Cells.Item(89, 13).Value = -18898.5
$ws.Cells.Item(89, 14).ClearContents()

$ws.Cells.Item(98, 8).Value = 385.85715
$ws.Cells.Item(98, 10).Value = 666.6667
$ws.Cells.Item(98, 12).Value = 2000.0001
$ws.Cells.Item(98, 14).Value = -4996.0001

$ws.Cells.Item(114, 8).Value = 1625.8
$ws.Cells.Item(114, 9).Value = 1164.8334
$ws.Cells.Item(114, 10).Value = 1933.1111
$ws.Cells.Item(114, 11).Value = 3494.5002
$ws.Cells.Item(114, 12).Value = 5799.3333
$ws.Cells.Item(114, 13).Value = -240.5001999999999
$ws.Cells.Item(114, 14).Value = -12307.3333

$ws.Cells.Item(117, 8).Value = 1093
$ws.Cells.Item(117, 9).Value = 253.33333
$ws.Cells.Item(117, 10).Value = 1932.6666
$ws.Cells.Item(117, 11).Value = 759.99999
$ws.Cells.Item(117, 12).Value = 5797.9998
$ws.Cells.Item(117, 13).Value = 2682.00001
$ws.Cells.Item(117, 14).Value = -12681.9998

$ws.Cells.Item(135, 8).Value = 986.2093
$ws.Cells.Item(135, 9).Value = 589.1613
$ws.Cells.Item(135, 10).Value = 2011.9166
$ws.Cells.Item(135, 11).Value = 5302.4517
$ws.Cells.Item(135, 12).Value = 18107.2494
$ws.Cells.Item(135, 13).Value = -2767.4517
$ws.Cells.Item(135, 14).Value = -23177.2494

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(17, 8).Value = 30000
$ws.Cells.Item(17, 10).Value = 10000
$ws.Cells.Item(17, 12).Value = 10000
$ws.Cells.Item(17, 14).Value = -10336

$ws.Cells.Item(113, 8).Value = 1970.72
$ws.Cells.Item(113, 9).Value = 1411
$ws.Cells.Item(113, 11).Value = 1411
$ws.Cells.Item(113, 13).Value = 759

$ws.Cells.Item(122, 8).Value = 696132.6
$ws.Cells.Item(122, 9).Value = 856201.0600000001
$ws.Cells.Item(122, 10).Value = 2502.6667
$ws.Cells.Item(122, 11).Value = 2568603.18
$ws.Cells.Item(122, 12).Value = 7508.000100000001
$ws.Cells.Item(122, 13).Value = -2566153.18
$ws.Cells.Item(122, 14).Value = -12408.0001

$ws.Cells.Item(137, 8).Value = 55000
$ws.Cells.Item(137, 10).Value = 55000
$ws.Cells.Item(137, 12).Value = 55000
$ws.Cells.Item(137, 14).Value = -65200

$ws.Cells.Item(138, 8).Value = 65466.668
$ws.Cells.Item(138, 10).Value = 65466.668
$ws.Cells.Item(138, 12).Value = 65466.668
$ws.Cells.Item(138, 14).Value = -75746.66800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(3, 8).Value = 5000000
$ws.Cells.Item(3, 10).Value = 5000000
$ws.Cells.Item(3, 12).Value = 5000000
$ws.Cells.Item(3, 14).Value = -5000224

$ws.Cells.Item(15, 8).Value = 5000000
$ws.Cells.Item(15, 10).Value = 5000000
$ws.Cells.Item(15, 12).Value = 5000000
$ws.Cells.Item(15, 14).Value = -5000340

$ws.Cells.Item(39, 8).Value = 26708
$ws.Cells.Item(39, 9).Value = 25029.5
$ws.Cells.Item(39, 11).Value = 25029.5
$ws.Cells.Item(39, 13).Value = -24569.5

